# Convert the colour-name lookup columns (AF:AQ) into PsychoPy-style
# normalized RGB triplets ([-1, 1] range) instead of CSS colour-name
# strings, replacing the old IFS(...) formulas with literal values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colorMap = @{
    "dodgerblue" = "[-0.843137254901961, -0.929411764705882, 0.568627450980392]"
    "yellow"     = "[0.984313725490196, 0.992156862745098, -0.835294117647059]"
    "limegreen"  = "[-1, 0.96078431372549, -0.976470588235294]"
    "aqua"       = "[-0.843137254901961, 0.984313725490196, 1]"
    "white"      = "[0.888, 0.888, 0.888]"
    "magenta"    = "[1, -0.67843137254902, 1]"
    "black"      = "[-0.23922, -0.41176, -0.34902]"
}

# Columns AF..AQ = 32..43, rows 2..4 hold the per-trial colour-name results.
for ($r = 2; $r -le 4; $r++) {
    for ($c = 32; $c -le 43; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $colorName = $cell.Value2
        $cell.Value = $colorMap[$colorName]
    }
}

# Widen the newly-converted columns so the long RGB-array text is visible,
# and move the viewport / selection over to the edited area.
$ws.Range("AF1:AQ4").Columns.AutoFit()
$ws.Range("AP4").Select()
$excel.ActiveWindow.Zoom = 100
